$wb = $excel.ActiveWorkbook

# Sheet: indexedListAsLeafTestLeft
$wsLeft = $wb.Worksheets.Item("indexedListAsLeafTestLeft")
$wsLeft.Range("C1").Value = "listAsLeafTestLeft[0]#test?readAs=text"
$wsLeft.Range("E1").Value = "listAsLeafTestLeft[0]#list[1]?readAs=text"
$wsLeft.Range("J1").Value = "listAsLeafTestLeft[1]#list[1]?readAs=text"

# Sheet: indexedListAsLeafTestOption
$wsOption = $wb.Worksheets.Item("indexedListAsLeafTestOption")
$wsOption.Range("C1").Value = "listAsLeafTestOption#test?readAs=text"
$wsOption.Range("E1").Value = "listAsLeafTestOption#list[1]?readAs=text"

# On the no-longer-active sheet, reset its selection to C1 (matches diff)
$wsLeft.Range("C1").Select()

# Change active sheet selection: indexedListAsLeaf becomes active
$wsLeaf = $wb.Worksheets.Item("indexedListAsLeaf")
$wsLeaf.Activate()
$wsLeaf.Range("E2").Select()

$wb.Save()
